$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44729
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 7000
$ws.Range("M2").Value = 7000
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 194

$ws.Range("D3").Value = 44690
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 7000
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 194

$ws.Range("D4").Value = 44372
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = 7000
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 194

$ws.Range("D5").Value = 44756
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 13000
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 361

$ws.Range("D6").Value = 44741
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 9000
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 9000
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 250

$ws.Range("D7").Value = 44750
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = 9000
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 250

$ws.Range("D8").Value = 44715
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 9000
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 9000
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 250

$ws.Range("D9").Value = 44719
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = 9000
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 250

$ws.Range("D10").Value = 44720
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 9000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 9000
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 250

$ws.Range("D11").Value = 44755
$ws.Range("J11").Value = 150
$ws.Range("K11").Value = 8500
$ws.Range("L11").Value = 8500
$ws.Range("M11").Value = 8500
$ws.Range("O11").Value = "Región Metropolitana"
$ws.Range("P11").Value = 236

$ws.Range("D12").Value = 44707
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 9000
$ws.Range("L12").Value = 9000
$ws.Range("M12").Value = 9000
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 250

$ws.Range("D13").Value = 44701
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 7000
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 194

$ws.Range("D14").Value = 44386
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 6500
$ws.Range("L14").Value = 6500
$ws.Range("M14").Value = 6500
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 181

$ws.Range("D15").Value = 44706
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 9000
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 250

$ws.Range("D16").Value = 44342
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 7000
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 194

$ws.Range("D17").Value = 44376
$ws.Range("J17").Value = 150
$ws.Range("K17").Value = 6500
$ws.Range("L17").Value = 6500
$ws.Range("M17").Value = 6500
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 181

$ws.Range("D18").Value = 44364
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 7000
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 194

$ws.Range("D19").Value = 44753
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 9000
$ws.Range("L19").Value = 9000
$ws.Range("M19").Value = 9000
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 250

$ws.Range("D20").Value = 44711
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 8500
$ws.Range("L20").Value = 8500
$ws.Range("M20").Value = 8500
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 236

$ws.Range("D21").Value = 44725
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 8000
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = 8000
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 222

$ws.Range("D22").Value = 44736
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 9000
$ws.Range("L22").Value = 9000
$ws.Range("M22").Value = 9000
$ws.Range("O22").Value = "Región Metropolitana"
$ws.Range("P22").Value = 250

$ws.Range("D23").Value = 44746
$ws.Range("J23").Value = 150
$ws.Range("K23").Value = 8000
$ws.Range("L23").Value = 8000
$ws.Range("M23").Value = 8000
$ws.Range("O23").Value = "Región del Maule"
$ws.Range("P23").Value = 222

$ws.Range("D24").Value = 44348
$ws.Range("J24").Value = 150
$ws.Range("K24").Value = 7000
$ws.Range("L24").Value = 7000
$ws.Range("M24").Value = 7000
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 194

$ws.Range("D25").Value = 44354
$ws.Range("J25").Value = 150
$ws.Range("K25").Value = 7000
$ws.Range("L25").Value = 7000
$ws.Range("M25").Value = 7000
$ws.Range("O25").Value = "Región del Maule"
$ws.Range("P25").Value = 194

$ws.Range("D26").Value = 44748
$ws.Range("J26").Value = 150
$ws.Range("K26").Value = 8000
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = 8000
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 222

$ws.Range("D27").Value = 44371
$ws.Range("J27").Value = 150
$ws.Range("K27").Value = 6500
$ws.Range("L27").Value = 6500
$ws.Range("M27").Value = 6500
$ws.Range("O27").Value = "Región Metropolitana"
$ws.Range("P27").Value = 181

$ws.Range("D28").Value = 44726
$ws.Range("J28").Value = 150
$ws.Range("K28").Value = 8000
$ws.Range("L28").Value = 8000
$ws.Range("M28").Value = 8000
$ws.Range("O28").Value = "Región del Maule"
$ws.Range("P28").Value = 222

$ws.Range("D29").Value = 44340
$ws.Range("J29").Value = 150
$ws.Range("K29").Value = 7000
$ws.Range("L29").Value = 7000
$ws.Range("M29").Value = 7000
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 194

$ws.Range("D30").Value = 44362
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 6500
$ws.Range("L30").Value = 6500
$ws.Range("M30").Value = 6500
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 181

$ws.Range("D31").Value = 44747
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 9000
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = 9000
$ws.Range("O31").Value = "Región Metropolitana"
$ws.Range("P31").Value = 250

$ws.Range("D32").Value = 44757
$ws.Range("J32").Value = 200
$ws.Range("K32").Value = 14000
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = 14000
$ws.Range("O32").Value = "Región Metropolitana"
$ws.Range("P32").Value = 389

$ws.Range("D33").Value = 44355
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 7000
$ws.Range("L33").Value = 7000
$ws.Range("M33").Value = 7000
$ws.Range("O33").Value = "Región Metropolitana"
$ws.Range("P33").Value = 194

$ws.Range("D34").Value = 44358
$ws.Range("J34").Value = 150
$ws.Range("K34").Value = 7000
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = 7000
$ws.Range("O34").Value = "Región Metropolitana"
$ws.Range("P34").Value = 194
